$wb = $excel.ActiveWorkbook

# ----- Sheet: 展览 -----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1382
$ws.Range("F3").Value = 1380
$ws.Range("F5").Value = 211
$ws.Range("F6").Value = 644
$ws.Range("F7").Value = 19
$ws.Range("F8").Value = 336
$ws.Range("F9").Value = 461
$ws.Range("F10").Value = 68
$ws.Range("F11").Value = 1347
$ws.Range("F12").Value = 31079
$ws.Range("F13").Value = 6556
$ws.Range("F14").Value = 87
$ws.Range("F15").Value = 327
$ws.Range("F16").Value = 550
$ws.Range("F17").Value = 95
$ws.Range("F19").Value = 73
$ws.Range("F20").Value = 37
$ws.Range("F21").Value = 410
$ws.Range("F22").Value = 83
$ws.Range("F23").Value = 752
$ws.Range("F24").Value = 305
$ws.Range("F25").Value = 362
$ws.Range("F26").Value = 410
$ws.Range("F28").Value = 156
$ws.Range("F29").Value = 36
$ws.Range("F30").Value = 712
$ws.Range("F31").Value = 261
$ws.Range("F32").Value = 127
$ws.Range("F33").Value = 677
$ws.Range("F34").Value = 96
$ws.Range("F36").Value = 758
$ws.Range("F37").Value = 269

# ----- Sheet: 演出 -----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 1093
$ws.Range("F5").Value = 125
$ws.Range("F6").Value = 288
$ws.Range("F7").Value = 4306
$ws.Range("F9").Value = 222
$ws.Range("F10").Value = 3
$ws.Range("F19").Value = 4270

# ----- Sheet: 本地生活 -----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("G2").Value = "不可售"
$ws.Range("F4").Value = 1377
$ws.Range("F5").Value = 330

# ----- Sheet: 全部类型 -----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1377
$ws.Range("F3").Value = 330
$ws.Range("F4").Value = 1093
$ws.Range("F5").Value = 1382
$ws.Range("F7").Value = 1380
$ws.Range("F8").Value = 211
$ws.Range("F9").Value = 644
$ws.Range("F10").Value = 19
$ws.Range("F11").Value = 336
$ws.Range("F12").Value = 461
$ws.Range("F14").Value = 68
$ws.Range("F15").Value = 1347
$ws.Range("F16").Value = 125
$ws.Range("F17").Value = 125
$ws.Range("F18").Value = 288
$ws.Range("F20").Value = 222
$ws.Range("F21").Value = 222
$ws.Range("F22").Value = 3
$ws.Range("F24").Value = 327
$ws.Range("F26").Value = 550
$ws.Range("F27").Value = 95
$ws.Range("F29").Value = 73
$ws.Range("F31").Value = 37
$ws.Range("F33").Value = 410
$ws.Range("F34").Value = 83
$ws.Range("F35").Value = 752
$ws.Range("F36").Value = 305
$ws.Range("F37").Value = 362
$ws.Range("F38").Value = 410
$ws.Range("F40").Value = 156
$ws.Range("F41").Value = 36
$ws.Range("F42").Value = 712
$ws.Range("F44").Value = 261
$ws.Range("F45").Value = 127
$ws.Range("F46").Value = 96
$ws.Range("F47").Value = 758
$ws.Range("F48").Value = 269
